{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// 1) Update the activation date 01/01/2020 -> 01/01/2023.\n//    Use a document-wide search/replace so formatting (the existing run) is preserved.\nconst searchResults = body.search(\"Ativa\u00e7\u00e3o: 01/01/2020\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"Ativa\u00e7\u00e3o: 01/01/2023\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Insert new italic English paragraphs right after three specific Portuguese\n//    paragraphs, matching the exact (trimmed) text of each.\nconst insertions = [\n  {\n    match: \"Apresentar conceitos sobre fontes renov\u00e1veis para gera\u00e7\u00e3o de energia t\u00e9rmica, el\u00e9trica e veicular, dentre outras.\",\n    text: \"To present concepts about renewable sources for the generation of thermal, electrical and vehicular energy, among others\",\n  },\n  {\n    match: \"Fontes renov\u00e1veis e tecnologias limpas para gera\u00e7\u00e3o de energia. Estudo dos sistemas atuais nacionais e mundiais.\",\n    text: \"Renewable sources and clean technologies for energy generation. Study of current national and global systems.\",\n  },\n  {\n    match: \"Sistemas energ\u00e9ticos nacionais e mundiais: fontes renov\u00e1veis e f\u00f3sseis. Gera\u00e7\u00e3o de energia por fontes renov\u00e1veis: solar t\u00e9rmica e fotovoltaica; e\u00f3lica; mar\u00edtima. Gera\u00e7\u00e3o de biomassa para fins energ\u00e9ticos. Gerenciamento de res\u00edduos s\u00f3lidos urbanos: recicl\u00e1veis e n\u00e3o recicl\u00e1veis; programas empresariais para log\u00edstica reversa; a quest\u00e3o dos pol\u00edmeros; reflorestamento; processamento do lixo \u00famido dom\u00e9stico. Integra\u00e7\u00e3o de fontes renov\u00e1veis para gera\u00e7\u00e3o de energia: ciclos t\u00e9rmicos h\u00edbridos\",\n    text: \"National and global energy systems: renewable and fossil sources. Energy generation from renewable sources: solar thermal and photovoltaic; wind; maritime. Generation of biomass for energy purposes. Management of urban solid waste: recyclable and non-recyclable; enterprise programs for reverse logistics; the issue of polymers; reforestation; processing of domestic wet waste. Integration of renewable sources for energy generation: hybrid thermal cycles\",\n  },\n];\n\nfor (const { match, text } of insertions) {\n  const target = paragraphs.items.find((p) => p.text.trim() === match);\n  if (!target) {\n    throw new Error(\"Could not find paragraph: \" + match);\n  }\n  const newPara = target.insertParagraph(text, Word.InsertLocation.after);\n  newPara.font.set({ italic: true });\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# 1) Update the activation date 01/01/2020 -> 01/01/2023.\n$find = $d.Content.Find\n[void]$find.Execute(\"Ativa\u00e7\u00e3o: 01/01/2020\", $false, $false, $false, $false, $false, $true, 1, $false, \"Ativa\u00e7\u00e3o: 01/01/2023\", 2)\n\n# 2) Insert new italic English paragraphs right after three specific Portuguese\n#    paragraphs, matching the exact text of each (trailing paragraph marks\n#    trimmed off before comparing).\n$insertions = @(\n    @{\n        Match = \"Apresentar conceitos sobre fontes renov\u00e1veis para gera\u00e7\u00e3o de energia t\u00e9rmica, el\u00e9trica e veicular, dentre outras.\"\n        Text  = \"To present concepts about renewable sources for the generation of thermal, electrical and vehicular energy, among others\"\n    },\n    @{\n        Match = \"Fontes renov\u00e1veis e tecnologias limpas para gera\u00e7\u00e3o de energia. Estudo dos sistemas atuais nacionais e mundiais.\"\n        Text  = \"Renewable sources and clean technologies for energy generation. Study of current national and global systems.\"\n    },\n    @{\n        Match = \"Sistemas energ\u00e9ticos nacionais e mundiais: fontes renov\u00e1veis e f\u00f3sseis. Gera\u00e7\u00e3o de energia por fontes renov\u00e1veis: solar t\u00e9rmica e fotovoltaica; e\u00f3lica; mar\u00edtima. Gera\u00e7\u00e3o de biomassa para fins energ\u00e9ticos. Gerenciamento de res\u00edduos s\u00f3lidos urbanos: recicl\u00e1veis e n\u00e3o recicl\u00e1veis; programas empresariais para log\u00edstica reversa; a quest\u00e3o dos pol\u00edmeros; reflorestamento; processamento do lixo \u00famido dom\u00e9stico. Integra\u00e7\u00e3o de fontes renov\u00e1veis para gera\u00e7\u00e3o de energia: ciclos t\u00e9rmicos h\u00edbridos\"\n        Text  = \"National and global energy systems: renewable and fossil sources. Energy generation from renewable sources: solar thermal and photovoltaic; wind; maritime. Generation of biomass for energy purposes. Management of urban solid waste: recyclable and non-recyclable; enterprise programs for reverse logistics; the issue of polymers; reforestation; processing of domestic wet waste. Integration of renewable sources for energy generation: hybrid thermal cycles\"\n    }\n)\n\nforeach ($item in $insertions) {\n    $target = $null\n    foreach ($p in $d.Paragraphs) {\n        $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($ptext -eq $item.Match) {\n            $target = $p\n            break\n        }\n    }\n    if ($target -eq $null) {\n        throw \"Could not find paragraph: \" + $item.Match\n    }\n\n    [void]$target.Range.InsertParagraphAfter()\n\n    $newRange = $target.Next().Range\n    [void]$newRange.MoveEnd(1, -1)\n    $newRange.Text = $item.Text\n    $newRange.Font.Italic = $true\n}\n"}
